# Update countries & provincias Spain
# Applies:
#  - Refresh timestamp string in A1
#  - Update totals for Estados Unidos (row 4) and Kazajistan (row 76)
#  - Insert "Honduras" (with updated stats) ahead of "Taiwan" in the
#    country ranking, shifting Taiwan/Ghana/Jordania/Reunion down by one row
#  - Insert "Fiyi" (with updated stats) ahead of "Laos" in the country
#    ranking, shifting Laos/Namibia/Mongolia/Dominica down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 10 de Abril de 2020 a las 05:22"

# --- Estados Unidos (row 4) updated totals -----------------------------
$ws.Cells.Item(4,2).Value = 468887
$ws.Cells.Item(4,3).Value = 321
$ws.Cells.Item(4,5).Value = 426262
$ws.Cells.Item(4,7).Value = 6
$ws.Cells.Item(4,8).Value = 16697

# --- Kazajistan (row 76) updated totals --------------------------------
$ws.Cells.Item(76,2).Value = 802
$ws.Cells.Item(76,3).Value = 21
$ws.Cells.Item(76,5).Value = 733
$ws.Cells.Item(76,7).Value = 1
$ws.Cells.Item(76,8).Value = 9

# --- Honduras inserted before Taiwan (rows 96-100) ---------------------
$ws.Cells.Item(96,1).Value = "Honduras"
$ws.Cells.Item(96,2).Value = 382
$ws.Cells.Item(96,3).Value = 39
$ws.Cells.Item(96,4).Value = 6
$ws.Cells.Item(96,5).Value = 353
$ws.Cells.Item(96,6).Value = 10
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 23

$ws.Cells.Item(97,1).Value = "Taiwan"
$ws.Cells.Item(97,2).Value = 380
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 80
$ws.Cells.Item(97,5).Value = 295
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 5

$ws.Cells.Item(98,1).Value = "Ghana"
$ws.Cells.Item(98,2).Value = 378
$ws.Cells.Item(98,3).Value = 0
$ws.Cells.Item(98,4).Value = 3
$ws.Cells.Item(98,5).Value = 369
$ws.Cells.Item(98,6).Value = 2
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 6

$ws.Cells.Item(99,1).Value = "Jordania"
$ws.Cells.Item(99,2).Value = 372
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 161
$ws.Cells.Item(99,5).Value = 204
$ws.Cells.Item(99,6).Value = 5
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 7

$ws.Cells.Item(100,1).Value = "Reunion"
$ws.Cells.Item(100,2).Value = 362
$ws.Cells.Item(100,3).Value = 0
$ws.Cells.Item(100,4).Value = 40
$ws.Cells.Item(100,5).Value = 322
$ws.Cells.Item(100,6).Value = 4
$ws.Cells.Item(100,7).Value = 0
$ws.Cells.Item(100,8).Value = 0

# --- Fiyi inserted before Laos (rows 172-176) ---------------------------
$ws.Cells.Item(172,1).Value = "Fiyi"
$ws.Cells.Item(172,2).Value = 16
$ws.Cells.Item(172,3).Value = 1
$ws.Cells.Item(172,4).Value = 0
$ws.Cells.Item(172,5).Value = 16
$ws.Cells.Item(172,6).Value = 0
$ws.Cells.Item(172,7).Value = 0
$ws.Cells.Item(172,8).Value = 0

$ws.Cells.Item(173,1).Value = "Laos"
$ws.Cells.Item(173,2).Value = 16
$ws.Cells.Item(173,3).Value = 0
$ws.Cells.Item(173,4).Value = 0
$ws.Cells.Item(173,5).Value = 16
$ws.Cells.Item(173,6).Value = 0
$ws.Cells.Item(173,7).Value = 0
$ws.Cells.Item(173,8).Value = 0

$ws.Cells.Item(174,1).Value = "Namibia"
$ws.Cells.Item(174,2).Value = 16
$ws.Cells.Item(174,3).Value = 0
$ws.Cells.Item(174,4).Value = 3
$ws.Cells.Item(174,5).Value = 13
$ws.Cells.Item(174,6).Value = 0
$ws.Cells.Item(174,7).Value = 0
$ws.Cells.Item(174,8).Value = 0

$ws.Cells.Item(175,1).Value = "Mongolia"
$ws.Cells.Item(175,2).Value = 16
$ws.Cells.Item(175,3).Value = 0
$ws.Cells.Item(175,4).Value = 4
$ws.Cells.Item(175,5).Value = 12
$ws.Cells.Item(175,6).Value = 0
$ws.Cells.Item(175,7).Value = 0
$ws.Cells.Item(175,8).Value = 0

$ws.Cells.Item(176,1).Value = "Dominica"
$ws.Cells.Item(176,2).Value = 16
$ws.Cells.Item(176,3).Value = 1
$ws.Cells.Item(176,4).Value = 5
$ws.Cells.Item(176,5).Value = 11
$ws.Cells.Item(176,6).Value = 0
$ws.Cells.Item(176,7).Value = 0
$ws.Cells.Item(176,8).Value = 0
